$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = '''69.632.91'
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = '''3.376.30'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +3.61%  '
$ws.Range("E4").Value = '  -0.01%  '
$c = $ws.Range("D5")
$c.Value = '''191.03'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.98%  '
$c = $ws.Range("D6")
$c.Value = '''593.20'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.97%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("E9").Value = '  +2.08%  '
$ws.Range("E10").Value = '  +2.40%  '
$c = $ws.Range("D11")
$c.Value = '''0.419'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.71%  '
$c = $ws.Range("D12")
$c.Value = '''3.965.80'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +3.76%  '
$ws.Range("E13").Value = '  -0.79%  '
$c = $ws.Range("D14")
$c.Value = '''28.61'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +3.39%  '
$c = $ws.Range("D15")
$c.Value = '''69.632.33'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +2.66%  '
$ws.Range("E16").Value = '  +1.77%  '
$c = $ws.Range("D17")
$c.Value = '''3.367.57'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +3.61%  '
$c = $ws.Range("D18")
$c.Value = '''454.13'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +15.27%  '
$c = $ws.Range("D20")
$c.Value = '''13.82'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.75%  '
$c = $ws.Range("D21")
$c.Value = '''7.77'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.27%  '
$c = $ws.Range("D22")
$c.Value = '''75.78'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +5.62%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +1.05%  '
$ws.Range("E25").Value = '  +3.26%  '
$ws.Range("E26").Value = '  +2.03%  '
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E29").Value = '  +3.21%  '
$c = $ws.Range("D30")
$c.Value = '''23.38'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +3.06%  '
$c = $ws.Range("D31")
$c.Value = '''5.59'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("E32").Value = '  +2.46%  '
$c = $ws.Range("D33")
$c.Value = '''7.00'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("E35").Value = '  +6.70%  '
$c = $ws.Range("D36")
$c.Value = '''164.59'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.04%  '
$ws.Range("E37").Value = '  +2.02%  '
$c = $ws.Range("D38")
$c.Value = '''27.75'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +4.08%  '
$c = $ws.Range("D39")
$c.Value = '''0.812'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.30%  '
$c = $ws.Range("D40")
$c.Value = '''4.60'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.28%  '
$c = $ws.Range("D41")
$c.Value = '''6.60'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.51%  '
$c = $ws.Range("D42")
$c.Value = '''2.739.64'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +4.95%  '
$ws.Range("E43").Value = '  +2.57%  '
$c = $ws.Range("D44")
$c.Value = '''25.46'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.65%  '
$c = $ws.Range("D45")
$c.Value = '''0.0688'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.09%  '
$c = $ws.Range("D46")
$c.Value = '''41.11'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.22%  '
$c = $ws.Range("D47")
$c.Value = '''340.46'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.99%  '
$ws.Range("E48").Value = '  +2.69%  '
$c = $ws.Range("D49")
$c.Value = '''32.87'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +6.63%  '
$ws.Range("E50").Value = '  +4.63%  '
$ws.Range("E51").Value = '  -0.72%  '
